$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1 (header): extend with P1=14, Q1=15, matching style of existing header cells ---
$ws.Range("O1").Copy() | Out-Null
$ws.Range("P1:Q1").PasteSpecial(-4122) | Out-Null
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# --- Rows 2-25: swap values in columns I,K,M,O (1<->2) and add P,Q = 2 ---
for ($r = 2; $r -le 25; $r++) {
    foreach ($col in @("I","K","M","O")) {
        $cell = $ws.Range("$col$r")
        if ($cell.Value2 -eq 1) {
            $cell.Value2 = 2
        } else {
            $cell.Value2 = 1
        }
    }
    $ws.Range("P$r").Value = 2
    $ws.Range("Q$r").Value = 2
}
